# Fix header labels on existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the header style (bold, bordered, centered) from the Weekly Quantity sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows
$data = @(
    @(44934.99999999999, 89, 56.36841223283109, 121.4190507364998),
    @(44948.99999999999, 69, 38.96756999976299, 100.8095716101138),
    @(44955.99999999999, 59, 26.30359341674797, 90.14447429606828),
    @(44969.99999999999, 39, 7.310635576944902, 69.51760775655991),
    @(44983.99999999999, 20, -11.43420082647092, 50.37346365444039),
    @(44990.99999999999, 10, -21.68093355876133, 41.79362947427179),
    @(44997.99999999999, 0, -33.67789143011462, 31.99726878261579),
    @(45004.99999999999, 0, -41.56842810817901, 24.54986481556175),
    @(45011.99999999999, 0, -50.66436001287508, 11.72109883429341),
    @(45018.99999999999, 0, -62.63988055626628, 2.083127133313158),
    @(45025.99999999999, 0, -69.85754186225162, -8.688016293594661),
    @(45032.99999999999, 0, -82.51546722852314, -17.70752555942549),
    @(45039.99999999999, 0, -92.46113933621206, -29.49301484399637)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Copy the date-format style from the Weekly Quantity sheet onto the ds column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A14").PasteSpecial(-4122)

$wsForecast.Range("A1").Select()
